$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.509.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.843.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.84"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5333"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3048"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06893"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.15"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07740"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7439"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.80%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.842.83"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.90"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.000"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.17%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.99"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.10%  "

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007947"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.528.00"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.087.86"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.624"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.991"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.313"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.32"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.207"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.689"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.84"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.270"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08788"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.060"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04801"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.935"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7282"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.31%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.108"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.309"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01719"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4780"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9145"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.46"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.885"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.486"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.41%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4134"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.99%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.040"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1243"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.05%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.87"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8979"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05799"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.67%  "
